$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Archivo Lote")
$ws.Range("A1").EntireColumn.Delete()
